$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-4 (A..J), replacing old rows 2-6.
# Column order: Hoje, Operadora, Data da Notificação, Demanda, Protocolo,
# Beneficiário, Prazo, Respondido, Natureza, Opções

$data = @(
    @("15-02-2023", "417823 - PREMIUM SAÚDE S.A", "14/02/2023  15:30:29", 12092359, 8514904, "MARIA EMILIA FONSECA RODRIGUES", "10 dias úteis", "NO", "Assistencial", "Responder  Detalhes"),
    @("15-02-2023", "417823 - PREMIUM SAÚDE S.A", "15/02/2023  09:41:19", 12093061, 8515747, "DALILA DE OLIVEIRA SILVA", "10 dias úteis", "NO", "No Assistencial", "Responder  Detalhes"),
    @("15-02-2023", "417823 - PREMIUM SAÚDE S.A", "15/02/2023  12:01:13", 12093493, 8516300, "VIVIANE KARINE SANTOS", "10 dias úteis", "NO", "Assistencial", "Responder  Detalhes")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

# Remove now-unused rows 5 and 6 (table shrank from 6 rows to 4 rows total).
$ws.Range("A5:J6").EntireRow.Delete()
